# Generate Report for Handoff
#
# The localization job for "5a56bfdc-9f51-4ef7-b32c-21cacbe62729.md" has
# moved from "In Translation" to "Ready for handoff": update its status,
# priority and timestamp cells across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 02:17:32"

# Status text got longer ("In Translation" -> "Ready for handoff"), so the
# zh-cn/de-de status columns widen to fit the new content.
$wsOverview.Columns.Item(5).ColumnWidth = 16.41
$wsOverview.Columns.Item(6).ColumnWidth = 16.41

# --- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-01 02:17:28"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.41

# --- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-01 02:17:32"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.41
